# Optima model-inputs.xlsx edit:
# 1. Update the "Model parameters" description for transnorm.
# 2. Add a new "format" column (H) to the "Data constants" sheet, classifying
#    each row as "percentage" or "number".
# 3. Update sheet view / active-tab state to match final navigation.

$wb = $excel.ActiveWorkbook

$wsModelParams = $wb.Worksheets.Item("Model parameters")
$wsDataInputs  = $wb.Worksheets.Item("Data inputs")
$wsDataConst   = $wb.Worksheets.Item("Data constants")

# --- 1. Model parameters: update the description text for "transnorm" ---
$wsModelParams.Range("B46").Value = "Normalization factor for transmissibility"

# --- 2. Data constants: add new column H ("format") ---
$wsDataConst.Columns.Item(8).ColumnWidth = 11.14

$wsDataConst.Range("H1").Value = "format"
$wsDataConst.Range("H1").Font.Bold = $true

# Rows grouped by subheading -- number-valued groups (written first so the
# "number" string is registered before "percentage" in the shared-string table)
$wsDataConst.Range("H9:H25").Value = "number"
$wsDataConst.Range("H53:H59").Value = "number"

# Rows grouped by subheading -- percentage-valued groups
$wsDataConst.Range("H2:H8").Value = "percentage"
$wsDataConst.Range("H26:H52").Value = "percentage"

# --- 3. Sheet view / navigation state ---
$wsModelParams.Activate()
$excel.ActiveWindow.FreezePanes = $false
$wsModelParams.Range("A35").Select()
$excel.ActiveWindow.FreezePanes = $true
$wsModelParams.Range("C46").Select()

$wsDataInputs.Activate()
$wsDataInputs.Range("A2").Select()

$wsDataConst.Activate()
$wsDataConst.Range("I12").Select()
